# Katalon AI: append new data row (row 3) to the "AI Generated" sheet and
# widen a few columns so the new text fits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments -------------------------------------------------
# Excel COM's Range/Columns.ColumnWidth is expressed in "characters"; the
# host stores width on a pixel grid with ~5/6 of a character of built-in
# padding baked in, so feed it (target - 5/6) to land exactly on the target
# stored width (29 / 18 / 17 / 16 characters respectively).
$ws.Columns.Item(1).ColumnWidth = 29 - 5/6   # A: 19 -> 29
$ws.Columns.Item(3).ColumnWidth = 18 - 5/6   # C: 17 -> 18
$ws.Columns.Item(4).ColumnWidth = 17 - 5/6   # D: 16 -> 17
$ws.Columns.Item(5).ColumnWidth = 16 - 5/6   # E: 15 -> 16

# --- New data row --------------------------------------------------------------
$ws.Range("A3").Value = "Adrress 25 f1 @#$%^&*!(#)#*"
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = "25 f1 first name"
$ws.Range("D3").Value = "25 f1 last name"
# These look numeric; prefix with an apostrophe so Excel stores them as text
# (matching the source data, which keeps phone/zip/state codes as strings).
$ws.Range("E3").Value = "'25012334567955"
$ws.Range("F3").Value = "25 F1 City"
$ws.Range("G3").Value = "'251"
$ws.Range("H3").Value = "'2501"
